# Update the "想去人数" (number of people wanting to go) column (F)
# on both the "展览" and "全部类型" worksheets, rows 2-14.

$wb = $excel.ActiveWorkbook

# row -> new value
$updates = @{
    2  = 41
    3  = 149
    4  = 59
    5  = 488
    6  = 1419
    7  = 689
    8  = 103
    9  = 189
    10 = 133
    11 = 183
    12 = 107
    13 = 158
    14 = 144
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
